$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.683.98'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '2.548.35'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.77'
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.78'
$ws.Range("E6").Value = '  +5.70%  '
$ws.Range("E7").Value = '  -1.06%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.11'
$ws.Range("E10").Value = '  +1.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0806'
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.42'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("D14").Value = '2.935.10'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.91'
$ws.Range("E15").Value = '  +6.11%  '
$ws.Range("D16").Value = '2.535.02'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.836'
$ws.Range("E17").Value = '  -1.60%  '
$ws.Range("D18").Value = '42.706.47'
$ws.Range("E18").Value = '  -0.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.82'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("E21").Value = '  -2.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.17'
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.47'
$ws.Range("E23").Value = '  -3.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.93'
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.53'
$ws.Range("E26").Value = '  -1.79%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '40.76'
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.36'
$ws.Range("E29").Value = '  -1.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.12'
$ws.Range("E30").Value = '  -2.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.11'
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("E32").Value = '  -3.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.77'
$ws.Range("E33").Value = '  +13.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0806'
$ws.Range("E34").Value = '  +1.00%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.63'
$ws.Range("E35").Value = '  -3.08%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.06'
$ws.Range("E36").Value = '  -1.94%  '
$ws.Range("E37").Value = '  -3.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.26'
$ws.Range("E38").Value = '  -6.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.112'
$ws.Range("E39").Value = '  -1.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.17'
$ws.Range("E41").Value = '  +9.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.23'
$ws.Range("E42").Value = '  +1.76%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.34'
$ws.Range("E43").Value = '  +2.14%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0299'
$ws.Range("E45").Value = '  -1.53%  '
$ws.Range("D46").Value = '1.976.44'
$ws.Range("E46").Value = '  -1.25%  '
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("D48").Value = '2.792.13'
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '81.31'
$ws.Range("E49").Value = '  -3.91%  '
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.40'
$ws.Range("E51").Value = '  -1.98%  '
